$wb = $excel.ActiveWorkbook

# --- Update "Training Data" sheet: column D (traffic_volume) values ---
$wsTrain = $wb.Worksheets.Item("Training Data")

$newValues = @{
    2 = 44447
    3 = 37837
    4 = 41237
    5 = 41959
    6 = 35553
    7 = 38344
    8 = 45383
    9 = 41827
    10 = 42682
    11 = 42976
    12 = 40622
    13 = 46961
    14 = 43589
    15 = 36325
    16 = 42946
    17 = 41084
    18 = 40889
    19 = 41440
    20 = 46369
    21 = 43689
    22 = 41268
    23 = 43659
    24 = 43198
    25 = 48228
    26 = 43363
    27 = 39855
    28 = 31954
    29 = 24986
    30 = 31734
    31 = 33890
    32 = 36918
    33 = 40151
    34 = 42158
    35 = 45120
    36 = 43154
    37 = 47257
    38 = 42234
    39 = 1338
    40 = 1138
    41 = 1250
    42 = 1326
    43 = 1393
    44 = 1426
    45 = 1434
    46 = 1393
    47 = 1450
    48 = 1396
    49 = 1590
    50 = 1232
    51 = 1272
    52 = 1343
    53 = 1325
    54 = 1259
    55 = 1274
    56 = 1440
    57 = 1346
    58 = 1408
    59 = 1465
    60 = 1324
    61 = 1544
}

foreach ($row in $newValues.Keys) {
    $wsTrain.Range("D$row").Value = $newValues[$row]
}

# --- Update "Testing Data" sheet: clear the large inline-string values in column D ---
$wsTest = $wb.Worksheets.Item("Testing Data")

for ($row = 2; $row -le 13; $row++) {
    $wsTest.Range("D$row").ClearContents()
}
